$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update filename in A2
$ws.Range("A2").Value = "Waves_005.txt"

# Update numeric measurement columns in row 2
$ws.Range("D2").Value = 147
$ws.Range("E2").Value = 28
$ws.Range("F2").Value = 436.61
$ws.Range("G2").Value = 10.69
$ws.Range("H2").Value = 20.28

$ws.Range("K2").Value = 4.67
$ws.Range("L2").Value = 0.6
$ws.Range("M2").Value = 0.01
$ws.Range("N2").Value = 14.38
$ws.Range("O2").Value = 1.02
$ws.Range("P2").Value = 0.02

$ws.Range("Q2").Value = 66.20999999999999
$ws.Range("R2").Value = 12.7
$ws.Range("S2").Value = 0.2
$ws.Range("T2").Value = 5.53

$ws.Range("W2").Value = 333.3
$ws.Range("X2").Value = 28.31
$ws.Range("Y2").Value = 0.44
$ws.Range("Z2").Value = 25.29

$ws.Range("AC2").Value = 19.07
$ws.Range("AD2").Value = 1.39
$ws.Range("AE2").Value = 0.02
$ws.Range("AF2").Value = 21.67
$ws.Range("AG2").Value = 1.76
$ws.Range("AH2").Value = 0.03
$ws.Range("AI2").Value = 24.45
$ws.Range("AJ2").Value = 0.74
$ws.Range("AK2").Value = 0.01
